# Generate Report for Handoff
# Adds two new localization entries (d24ca8b8-... and f435da6e-...) to the
# Overview / zh-cn / de-de tables, each growing from 2 data rows to 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview", columns A:G)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")

$row = $loOverview.ListRows.Add()
$r = $row.Range.Row
$wsOverview.Cells.Item($r, 1).Value = "d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md"
$wsOverview.Cells.Item($r, 3).Value = ".md"
$wsOverview.Cells.Item($r, 4).Value = ""
$wsOverview.Cells.Item($r, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 7).Value = "2016-08-13 02:49:09"
$wsOverview.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/d5e6d7d6a3e6c0f4b6a5c4d3e2f1a0b9c8d7e6f5/e2e/d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md", "", "", "e2e\d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md") | Out-Null

$row = $loOverview.ListRows.Add()
$r = $row.Range.Row
$wsOverview.Cells.Item($r, 1).Value = "f435da6e-8620-4ced-ada5-010da4e88e2b.md"
$wsOverview.Cells.Item($r, 3).Value = ".md"
$wsOverview.Cells.Item($r, 4).Value = ""
$wsOverview.Cells.Item($r, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($r, 7).Value = "2016-08-13 02:49:09"
$wsOverview.Cells.Item($r, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($r, 2), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4/e2e/f435da6e-8620-4ced-ada5-010da4e88e2b.md", "", "", "e2e\f435da6e-8620-4ced-ada5-010da4e88e2b.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn" / displayName "zh_cn", columns A:P)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")

$row = $loZhCn.ListRows.Add()
$r = $row.Range.Row
$wsZhCn.Cells.Item($r, 2).Value = ".md"
$wsZhCn.Cells.Item($r, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item($r, 4).Value = "e2e"
$wsZhCn.Cells.Item($r, 5).Value = "ht"
$wsZhCn.Cells.Item($r, 6).Value = "False"
$wsZhCn.Cells.Item($r, 7).Value = "d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.4c3e84e86f46b6c88c77c8f36d15ce738f4eb993.zh-cn.xlf"
$wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 02:49:00"
$wsZhCn.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($r, 9).Value = ""
$wsZhCn.Cells.Item($r, 10).Value = ""
$wsZhCn.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($r, 12).Value = ""
$wsZhCn.Cells.Item($r, 13).Value = "True"
$wsZhCn.Cells.Item($r, 14).Value = ""
$wsZhCn.Cells.Item($r, 15).Value = "False"
$wsZhCn.Cells.Item($r, 16).Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/d5e6d7d6a3e6c0f4b6a5c4d3e2f1a0b9c8d7e6f5/e2e/d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md", "", "", "d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md") | Out-Null

$row = $loZhCn.ListRows.Add()
$r = $row.Range.Row
$wsZhCn.Cells.Item($r, 2).Value = ".md"
$wsZhCn.Cells.Item($r, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item($r, 4).Value = "e2e"
$wsZhCn.Cells.Item($r, 5).Value = "ht"
$wsZhCn.Cells.Item($r, 6).Value = "False"
$wsZhCn.Cells.Item($r, 7).Value = "f435da6e-8620-4ced-ada5-010da4e88e2b.b3088bb14c0678feb12ca23ebcc107f047002d2b.zh-cn.xlf"
$wsZhCn.Cells.Item($r, 8).Value = "2016-08-13 02:49:00"
$wsZhCn.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($r, 9).Value = ""
$wsZhCn.Cells.Item($r, 10).Value = ""
$wsZhCn.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item($r, 12).Value = ""
$wsZhCn.Cells.Item($r, 13).Value = "True"
$wsZhCn.Cells.Item($r, 14).Value = ""
$wsZhCn.Cells.Item($r, 15).Value = "False"
$wsZhCn.Cells.Item($r, 16).Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4/e2e/f435da6e-8620-4ced-ada5-010da4e88e2b.md", "", "", "f435da6e-8620-4ced-ada5-010da4e88e2b.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de-de" / displayName "de_de", columns A:P)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")

$row = $loDeDe.ListRows.Add()
$r = $row.Range.Row
$wsDeDe.Cells.Item($r, 2).Value = ".md"
$wsDeDe.Cells.Item($r, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item($r, 4).Value = "e2e"
$wsDeDe.Cells.Item($r, 5).Value = "ht"
$wsDeDe.Cells.Item($r, 6).Value = "False"
$wsDeDe.Cells.Item($r, 7).Value = "d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.4c3e84e86f46b6c88c77c8f36d15ce738f4eb993.de-de.xlf"
$wsDeDe.Cells.Item($r, 8).Value = "2016-08-13 02:49:09"
$wsDeDe.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($r, 9).Value = ""
$wsDeDe.Cells.Item($r, 10).Value = ""
$wsDeDe.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($r, 12).Value = ""
$wsDeDe.Cells.Item($r, 13).Value = "True"
$wsDeDe.Cells.Item($r, 14).Value = ""
$wsDeDe.Cells.Item($r, 15).Value = "False"
$wsDeDe.Cells.Item($r, 16).Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/d5e6d7d6a3e6c0f4b6a5c4d3e2f1a0b9c8d7e6f5/e2e/d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md", "", "", "d24ca8b8-4a6a-4ea5-84af-2e4d8682bb6d.md") | Out-Null

$row = $loDeDe.ListRows.Add()
$r = $row.Range.Row
$wsDeDe.Cells.Item($r, 2).Value = ".md"
$wsDeDe.Cells.Item($r, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item($r, 4).Value = "e2e"
$wsDeDe.Cells.Item($r, 5).Value = "ht"
$wsDeDe.Cells.Item($r, 6).Value = "False"
$wsDeDe.Cells.Item($r, 7).Value = "f435da6e-8620-4ced-ada5-010da4e88e2b.b3088bb14c0678feb12ca23ebcc107f047002d2b.de-de.xlf"
$wsDeDe.Cells.Item($r, 8).Value = "2016-08-13 02:49:09"
$wsDeDe.Cells.Item($r, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($r, 9).Value = ""
$wsDeDe.Cells.Item($r, 10).Value = ""
$wsDeDe.Cells.Item($r, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item($r, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item($r, 12).Value = ""
$wsDeDe.Cells.Item($r, 13).Value = "True"
$wsDeDe.Cells.Item($r, 14).Value = ""
$wsDeDe.Cells.Item($r, 15).Value = "False"
$wsDeDe.Cells.Item($r, 16).Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item($r, 1), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4/e2e/f435da6e-8620-4ced-ada5-010da4e88e2b.md", "", "", "f435da6e-8620-4ced-ada5-010da4e88e2b.md") | Out-Null
